# adding filter, update search, update relasi kat_soal, sorting asc desc
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data rows (update existing rows 2-3 first: cekfix-* -> peserta-*) ----
# Row 2: peserta-1
$ws.Range("A2").Value = "peserta-1"
$ws.Range("B2").Value = "password123"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 11111111
$ws.Range("F2").Value = "peserta-1"

# Row 3: peserta-2
$ws.Range("A3").Value = "peserta-2"
$ws.Range("B3").Value = "password123"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 22222222
$ws.Range("F3").Value = "peserta-2"

# --- Header row (row 1) -----------------------------------------------
# Add a new "filter" header in G1, then backfill the filter values for
# the rows that already existed.
$ws.Range("G1").Value = "filter"
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2

# Row 4: peserta-3 (new)
$ws.Range("A4").Value = "peserta-3"
$ws.Range("B4").Value = "password123"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 33333333
$ws.Range("F4").Value = "peserta-3"
$ws.Range("G4").Value = 1

# Row 5: peserta-4 (new)
$ws.Range("A5").Value = "peserta-4"
$ws.Range("B5").Value = "password123"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 44444444
$ws.Range("F5").Value = "peserta-4"
$ws.Range("G5").Value = 2

# Row 6: peserta-5 (new)
$ws.Range("A6").Value = "peserta-5"
$ws.Range("B6").Value = "password123"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 55555555
$ws.Range("F6").Value = "peserta-5"
$ws.Range("G6").Value = 1

# --- Column widths (best effort; engine quantizes to 1/6-char steps) ------
$ws.Columns.Item(3).ColumnWidth = 8.6
$ws.Columns.Item(5).ColumnWidth = 11.6

# --- Phonetic info on the whole used range (noConversion) -----------------
$ws.Range("A1:G6").SetPhonetic()

# --- Selection / active cell ----------------------------------------------
$ws.Range("H7").Select()
